# Applies the "Created one table with main Rsq values for g/lmers" edit:
#  - model-formula paragraphs: "Julian_first_follicle^3" -> "Julian_first_follicle - 170"
#    (and, for the two paragraphs without "Year", a whitespace re-wrap)
#  - chi2 / p-value numbers recomputed for every data row of all four tables
#  - two of the four tables get a narrower 3rd (chi2) column, and one of
#    them also gets a wider 4th (p) column
#  - a handful of row heights shrink slightly (content got shorter)
#  - a few p-values cross the significance boundary, so their Bold flag flips

$d = $word.ActiveDocument

function Replace-InRange($range, [string]$old, [string]$new) {
    # Re-anchor via $d.Range(start, end): a Range object obtained straight
    # off Cell.Range does not keep its own Find cursor pinned to its start
    # in this host, so .Find.Execute on it silently restarts the search
    # from the top of the document. Rebuilding an equivalent Range with
    # Document.Range(start, end) fixes the Find cursor to that span.
    $scoped = $d.Range($range.Start, $range.End)
    # Replace = 1 (wdReplaceOne): only replace the single match inside the
    # supplied range instead of sweeping the whole story (wdReplaceAll
    # ignores the caller's range boundaries and would clobber identical
    # text - e.g. "<0.001***" - that repeats in other table cells).
    $null = $scoped.Find.Execute($old, $false, $false, $false, $false, $false, $true, 0, $false, $new, 1)
}

function Set-CellText($table, [int]$row, [int]$col, [string]$old, [string]$new, $bold) {
    $cell = $table.Cell($row, $col)
    Replace-InRange $cell.Range $old $new
    if ($null -ne $bold) {
        $table.Cell($row, $col).Range.Bold = $bold
    }
}

function Set-RowHeight($table, [int]$anchorRow, [int]$anchorCol, [double]$points) {
    $table.Cell($anchorRow, $anchorCol).Row.Height = $points
}

# ---------------------------------------------------------------------
# Model formula paragraphs
# ---------------------------------------------------------------------
Replace-InRange $d.Content `
    "Model: Julian_first_follicle^3 ~ Block + Year + (1 | Population/Family) +     City_dist + Transect_ID + City_dist:Transect_ID" `
    "Model: Julian_first_follicle - 170 ~ Block + Year + (1 | Population/Family) +     City_dist + Transect_ID + City_dist:Transect_ID"

Replace-InRange $d.Content `
    "Model: Julian_first_follicle^3 ~ Block + (1 | Population/Family) + City_dist +     Transect_ID + City_dist:Transect_ID" `
    "Model: Julian_first_follicle - 170 ~ Block + (1 | Population/Family) +     City_dist + Transect_ID + City_dist:Transect_ID"

Replace-InRange $d.Content `
    "Model: Julian_first_follicle^3 ~ Block + Year + (1 | Population/Family) +     Urb_score + Transect_ID + Urb_score:Transect_ID" `
    "Model: Julian_first_follicle - 170 ~ Block + Year + (1 | Population/Family) +     Urb_score + Transect_ID + Urb_score:Transect_ID"

Replace-InRange $d.Content `
    "Model: Julian_first_follicle^3 ~ Block + (1 | Population/Family) + Urb_score +     Transect_ID + Urb_score:Transect_ID" `
    "Model: Julian_first_follicle - 170 ~ Block + (1 | Population/Family) +     Urb_score + Transect_ID + Urb_score:Transect_ID"

# ---------------------------------------------------------------------
# Table 1: City_dist, all years of data
# ---------------------------------------------------------------------
$t1 = $d.Tables(1)
$t1.Columns(3).Width = 54.2   # 1818 dxa -> 1084 dxa

Set-RowHeight $t1 2 2 28.65   # Block row:  600 -> 573
Set-RowHeight $t1 3 2 28.7    # Year row:   600 -> 574

Set-CellText $t1 2 3 "157,511.302"   "5.873"   $false
Set-CellText $t1 2 4 "<0.001***"     "0.118"   $false

Set-CellText $t1 3 3 "3,595,370.318" "52.576"  $false
# Year's p-value stays "<0.001***" (bold) - unchanged

Set-CellText $t1 4 3 "5.774"  "6.469"  $false
Set-CellText $t1 4 4 "0.016*" "0.011*" $true

Set-CellText $t1 5 3 "5.541"  "1.906" $false
Set-CellText $t1 5 4 "0.019*" "0.167" $false

Set-CellText $t1 6 3 "0.605" "0.074" $false
Set-CellText $t1 6 4 "0.437" "0.785" $false

# ---------------------------------------------------------------------
# Table 2: City_dist, one year of data
# ---------------------------------------------------------------------
$t2 = $d.Tables(2)
$t2.Columns(3).Width = 54.2   # 1634 dxa -> 1084 dxa
$t2.Columns(4).Width = 56.6   # 1347 dxa -> 1132 dxa

Set-RowHeight $t2 2 2 28.7    # Block row: 600 -> 574

Set-CellText $t2 2 3 "264,101.558" "12.258"  $false
Set-CellText $t2 2 4 "<0.001***"   "0.007**" $true

Set-CellText $t2 3 3 "2.566" "2.754" $false
Set-CellText $t2 3 4 "0.109" "0.097" $false

Set-CellText $t2 4 3 "1.639" "0.208" $false
Set-CellText $t2 4 4 "0.201" "0.649" $false

Set-CellText $t2 5 3 "0.262" "0.000" $false
Set-CellText $t2 5 4 "0.609" "0.984" $false

# ---------------------------------------------------------------------
# Table 3: Urb_score, all years of data
# ---------------------------------------------------------------------
$t3 = $d.Tables(3)
$t3.Columns(3).Width = 69.45  # 1818 dxa -> 1389 dxa

Set-CellText $t3 2 3 "318,061.786" "4,472.739" $false
# (Intercept) p-value stays "<0.001***" (bold) - unchanged

Set-RowHeight $t3 3 2 28.65   # Block row: 600 -> 573
Set-RowHeight $t3 4 2 28.7    # Year row:  600 -> 574

Set-CellText $t3 3 3 "157,510.402" "6.568" $false
Set-CellText $t3 3 4 "<0.001***"   "0.087" $false

Set-CellText $t3 4 3 "3,595,370.564" "55.639" $false
# Year's p-value stays "<0.001***" (bold) - unchanged

Set-CellText $t3 5 3 "0.785" "0.755" $false
Set-CellText $t3 5 4 "0.376" "0.385" $false

Set-CellText $t3 6 3 "0.009" "0.612" $false
Set-CellText $t3 6 4 "0.926" "0.434" $false

Set-CellText $t3 7 3 "3.174" "4.228" $false
Set-CellText $t3 7 4 "0.075" "0.04*" $true

# ---------------------------------------------------------------------
# Table 4: Urb_score, one year of data
# ---------------------------------------------------------------------
$t4 = $d.Tables(4)
$t4.Columns(3).Width = 69.45  # 1634 dxa -> 1389 dxa

Set-CellText $t4 2 3 "290,635.516" "4,602.175" $false
# (Intercept) p-value stays "<0.001***" (bold) - unchanged

Set-RowHeight $t4 3 2 28.7    # Block row: 600 -> 574

Set-CellText $t4 3 3 "264,100.422" "11.743"  $false
Set-CellText $t4 3 4 "<0.001***"   "0.008**" $true

Set-CellText $t4 4 3 "2.850" "2.209" $false
Set-CellText $t4 4 4 "0.091" "0.137" $false

Set-CellText $t4 5 3 "0.990" "2.070" $false
Set-CellText $t4 5 4 "0.32"  "0.15"  $false

Set-CellText $t4 6 3 "4.834" "4.435" $false
Set-CellText $t4 6 4 "0.028*" "0.035*" $true

Write-Output "done"
